$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Add the new "Wetland_08" (time 16:00) rows at 89-97.
#    Doing this BEFORE re-labelling rows 86-88 keeps the existing shared
#    string "Wetland_08" alive (still referenced), so later edits append
#    new strings instead of overwriting/recycling this slot.
# ---------------------------------------------------------------------------
$data08 = New-Object 'object[,]' 9,4
$t08 = 0.66666666666666663

$data08[0,0] = "Wetland_08"; $data08[0,1] = $t08; $data08[0,2] = 1; $data08[0,3] = 61
$data08[1,0] = "Wetland_08"; $data08[1,1] = $t08; $data08[1,2] = 1; $data08[1,3] = 63
$data08[2,0] = "Wetland_08"; $data08[2,1] = $t08; $data08[2,2] = 1; $data08[2,3] = 17
$data08[3,0] = "Wetland_08"; $data08[3,1] = $t08; $data08[3,2] = 2; $data08[3,3] = 16
$data08[4,0] = "Wetland_08"; $data08[4,1] = $t08; $data08[4,2] = 2; $data08[4,3] = 66
$data08[5,0] = "Wetland_08"; $data08[5,1] = $t08; $data08[5,2] = 2; $data08[5,3] = 30
$data08[6,0] = "Wetland_08"; $data08[6,1] = $t08; $data08[6,2] = 3; $data08[6,3] = 35
$data08[7,0] = "Wetland_08"; $data08[7,1] = $t08; $data08[7,2] = 3; $data08[7,3] = 52
$data08[8,0] = "Wetland_08"; $data08[8,1] = $t08; $data08[8,2] = 3; $data08[8,3] = 9

$ws.Range("A89:D97").Value = $data08
$ws.Range("B89:B97").NumberFormat = "h:mm"

# ---------------------------------------------------------------------------
# 2. Re-label rows 86-88 from "Wetland_08" to "Wetland_11" (new string,
#    appended to the shared strings table after "Wetland_08").
# ---------------------------------------------------------------------------
$ws.Range("A86:A88").Value = "Wetland_11"

# ---------------------------------------------------------------------------
# 3. Add the new "Wetland_10" (time 16:45) rows at 98-106 (new string,
#    appended to the shared strings table after "Wetland_11").
# ---------------------------------------------------------------------------
$data10 = New-Object 'object[,]' 9,4
$t10 = 0.69791666666666663

$data10[0,0] = "Wetland_10"; $data10[0,1] = $t10; $data10[0,2] = 1; $data10[0,3] = 40
$data10[1,0] = "Wetland_10"; $data10[1,1] = $t10; $data10[1,2] = 1; $data10[1,3] = 45
$data10[2,0] = "Wetland_10"; $data10[2,1] = $t10; $data10[2,2] = 1; $data10[2,3] = 51
$data10[3,0] = "Wetland_10"; $data10[3,1] = $t10; $data10[3,2] = 2; $data10[3,3] = 34
$data10[4,0] = "Wetland_10"; $data10[4,1] = $t10; $data10[4,2] = 2; $data10[4,3] = 48
$data10[5,0] = "Wetland_10"; $data10[5,1] = $t10; $data10[5,2] = 2; $data10[5,3] = 27
$data10[6,0] = "Wetland_10"; $data10[6,1] = $t10; $data10[6,2] = 3; $data10[6,3] = 14
$data10[7,0] = "Wetland_10"; $data10[7,1] = $t10; $data10[7,2] = 3; $data10[7,3] = 20
$data10[8,0] = "Wetland_10"; $data10[8,1] = $t10; $data10[8,2] = 3; $data10[8,3] = 20

$ws.Range("A98:D106").Value = $data10
$ws.Range("B98:B106").NumberFormat = "h:mm"

# ---------------------------------------------------------------------------
# 4. Update the sheet view: scroll so row 73 is the top row, and select A98.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 73
$ws.Range("A98").Select()
